$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The shape "TextBox 15" holds the "For more detail..." / "please visit
# https://github.com/..." call-to-action text at the bottom of the poster.
$shp = $s.Shapes.Item("TextBox 15")
$tr = $shp.TextFrame.TextRange

# Second paragraph: "please visit https://github.com/KellyK81/data-mining"
# becomes "please visit https://github.com/KellyK81/deep-learning", typed
# back in as four separate runs (mirrors how the author's edit/spell-check
# pass re-split the sentence in the source file).
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = ""
[void]$para2.InsertAfter("please visit https://")
[void]$para2.InsertAfter("github.com")
[void]$para2.InsertAfter("/KellyK81/d")
[void]$para2.InsertAfter("eep-learning")
